# Update the threshold table on Sheet1:
#   B2 (alpha_distance_range, Min): 5.6 -> 5.5
#   C2 (alpha_distance_range, Max): 11.7 -> 11.5
#   B3 (beta_distance_range,  Min): 5.7 -> 5.5
#   C3 (beta_distance_range,  Max): 10.4 -> 10.5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 10.5
